$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.038642423396847
$arrBF[0,2] = 1.044516097984917
$arrBF[0,3] = 1.046499221168529
$arrBF[0,4] = 1.055492008821156
$ws.Range("B2:F2").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036868983144936
$arrIM[0,1] = 1.043738752650487
$arrIM[0,2] = 1.047286928701408
$arrIM[0,3] = 1.049264485977201
$arrIM[0,4] = 1.058232318552412
$ws.Range("I2:M2").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.039616177363465
$arrBF[0,2] = 1.04525444585178
$arrBF[0,3] = 1.047349578577586
$arrBF[0,4] = 1.056392440981088
$ws.Range("B3:F3").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.037043910859589
$arrIM[0,1] = 1.044357262356669
$arrIM[0,2] = 1.047836720561548
$arrIM[0,3] = 1.049926397693937
$arrIM[0,4] = 1.058945980818656
$ws.Range("I3:M3").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.040246583607186
$arrBF[0,2] = 1.04573209896655
$arrBF[0,3] = 1.047900433007013
$arrBF[0,4] = 1.056975590338033
$ws.Range("B4:F4").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.037155301556384
$arrIM[0,1] = 1.04475719213287
$arrIM[0,2] = 1.048191698784552
$arrIM[0,3] = 1.050354670623625
$arrIM[0,4] = 1.059407642136563
$ws.Range("I4:M4").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.04051168265867
$arrBF[0,2] = 1.045932877220013
$arrBF[0,3] = 1.048132158427587
$arrBF[0,4] = 1.057220866905731
$ws.Range("B5:F5").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.037201698745561
$arrIM[0,1] = 1.044925252833738
$arrIM[0,2] = 1.048340745256467
$arrIM[0,3] = 1.050534708735926
$arrIM[0,4] = 1.059601693209535
$ws.Range("I5:M5").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.040556198402184
$arrBF[0,2] = 1.045966587140691
$arrBF[0,3] = 1.048171074680907
$arrBF[0,4] = 1.057262056950499
$ws.Range("B6:F6").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.037209463715757
$arrIM[0,1] = 1.044953466844553
$arrIM[0,2] = 1.048365759853681
$arrIM[0,3] = 1.050564937436922
$arrIM[0,4] = 1.059634273365299
$ws.Range("I6:M6").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.040250125575135
$arrBF[0,2] = 1.045734781883747
$arrBF[0,3] = 1.047903528760565
$arrBF[0,4] = 1.056978867263522
$ws.Range("B7:F7").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.037155923214151
$arrIM[0,1] = 1.044759438043598
$arrIM[0,2] = 1.048193691084524
$arrIM[0,3] = 1.050357076333805
$arrIM[0,4] = 1.059410235181189
$ws.Range("I7:M7").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.038971440833103
$arrBF[0,2] = 1.044765647473418
$arrBF[0,3] = 1.046786475023407
$arrBF[0,4] = 1.055796207737536
$ws.Range("B8:F8").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036928472896355
$arrIM[0,1] = 1.04394783981248
$arrIM[0,2] = 1.047472893007043
$arrIM[0,3] = 1.04948818733812
$arrIM[0,4] = 1.058473529398671
$ws.Range("I8:M8").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.036720742773146
$arrBF[0,2] = 1.043057149233229
$arrBF[0,3] = 1.044822862425213
$arrBF[0,4] = 1.053716174823376
$ws.Range("B9:F9").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.03651392792973
$arrIM[0,1] = 1.042515542268365
$arrIM[0,2] = 1.04619687624649
$arrIM[0,3] = 1.047956924980987
$arrIM[0,4] = 1.056822022952031
$ws.Range("I9:M9").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035222014608787
$arrBF[0,2] = 1.041917722671085
$arrBF[0,3] = 1.043517078096145
$arrBF[0,4] = 1.052332233614569
$ws.Range("B10:F10").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036228362666383
$arrIM[0,1] = 1.041559283414111
$arrIM[0,2] = 1.045342310310746
$arrIM[0,3] = 1.046936032291316
$arrIM[0,4] = 1.055720476266042
$ws.Range("I10:M10").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034573470977231
$arrBF[0,2] = 1.041424254068164
$arrBF[0,3] = 1.042952455624236
$arrBF[0,4] = 1.05173364033933
$ws.Range("B11:F11").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036102535716332
$arrIM[0,1] = 1.041144893814154
$arrIM[0,2] = 1.044971364806084
$arrIM[0,3] = 1.046493975189003
$arrIM[0,4] = 1.055243380175141
$ws.Range("I11:M11").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034332636175689
$arrBF[0,2] = 1.041240945767739
$arrBF[0,3] = 1.042742849758554
$arrBF[0,4] = 1.051511396983049
$ws.Range("B12:F12").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036055471810863
$arrIM[0,1] = 1.04099092339964
$arrIM[0,2] = 1.044833443002753
$arrIM[0,3] = 1.046329775916221
$arrIM[0,4] = 1.055066148667759
$ws.Range("I12:M12").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034384293187135
$arrBF[0,2] = 1.041280266546077
$arrBF[0,3] = 1.042787805475136
$arrBF[0,4] = 1.051559064353082
$ws.Range("B13:F13").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036065581928137
$arrIM[0,1] = 1.041023952718935
$arrIM[0,2] = 1.04486303384154
$arrIM[0,3] = 1.046364997182809
$arrIM[0,4] = 1.055104166168982
$ws.Range("I13:M13").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034553562204236
$arrBF[0,2] = 1.041409101985769
$arrBF[0,3] = 1.042935127081541
$arrBF[0,4] = 1.051715267569474
$ws.Range("B14:F14").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036098652049374
$arrIM[0,1] = 1.041132167538892
$arrIM[0,2] = 1.04495996691789
$arrIM[0,3] = 1.046480402411642
$arrIM[0,4] = 1.055228730499238
$ws.Range("I14:M14").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.03465786282654
$arrBF[0,2] = 1.041488480190231
$arrBF[0,3] = 1.043025912732169
$arrBF[0,4] = 1.051811522927956
$ws.Range("B15:F15").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036118984440454
$arrIM[0,1] = 1.041198835971372
$arrIM[0,2] = 1.045019672577578
$arrIM[0,3] = 1.046551507466972
$arrIM[0,4] = 1.055305476504573
$ws.Range("I15:M15").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035265065132704
$arrBF[0,2] = 1.041950470788316
$arrBF[0,3] = 1.043554566997882
$arrBF[0,4] = 1.052371974347639
$ws.Range("B16:F16").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036236667610283
$arrIM[0,1] = 1.041586778353697
$arrIM[0,2] = 1.04536690959687
$arrIM[0,3] = 1.046965370169448
$arrIM[0,4] = 1.055752137148386
$ws.Range("I16:M16").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035646058880331
$arrBF[0,2] = 1.042240242301939
$arrBF[0,3] = 1.043886390449608
$arrBF[0,4] = 1.052723708993158
$ws.Range("B17:F17").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036309905216252
$arrIM[0,1] = 1.041830038459855
$arrIM[0,2] = 1.045584478535728
$arrIM[0,3] = 1.047224974985122
$arrIM[0,4] = 1.056032284471852
$ws.Range("I17:M17").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035868326132189
$arrBF[0,2] = 1.042409252409721
$arrBF[0,3] = 1.044080013667689
$arrBF[0,4] = 1.052928933684776
$ws.Range("B18:F18").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036352413544381
$arrIM[0,1] = 1.041971896719845
$arrIM[0,2] = 1.045711294675645
$arrIM[0,3] = 1.047376397645896
$arrIM[0,4] = 1.056195678099087
$ws.Range("I18:M18").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035944120290729
$arrBF[0,2] = 1.04246687897006
$arrBF[0,3] = 1.044146047070064
$arrBF[0,4] = 1.052998920786525
$ws.Range("B19:F19").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.0363668721587
$arrIM[0,1] = 1.042020261393813
$arrIM[0,2] = 1.045754520720961
$arrIM[0,3] = 1.047428028755885
$arrIM[0,4] = 1.056251389107502
$ws.Range("I19:M19").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035605177709412
$arrBF[0,2] = 1.042209153452461
$arrBF[0,3] = 1.04385078104142
$arrBF[0,4] = 1.052685964583299
$ws.Range("B20:F20").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036302069223453
$arrIM[0,1] = 1.041803942177913
$arrIM[0,2] = 1.045561144557122
$arrIM[0,3] = 1.047197121883161
$arrIM[0,4] = 1.056002228503619
$ws.Range("I20:M20").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034503714947546
$arrBF[0,2] = 1.041371163482926
$arrBF[0,3] = 1.042891741216017
$arrBF[0,4] = 1.051669266821289
$ws.Range("B21:F21").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.03608892272178
$arrIM[0,1] = 1.041100302280277
$arrIM[0,2] = 1.04493142629254
$arrIM[0,3] = 1.04644641842866
$arrIM[0,4] = 1.055192049862816
$ws.Range("I21:M21").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.033811547500661
$arrBF[0,2] = 1.040844216641193
$arrBF[0,3] = 1.042289450064876
$arrBF[0,4] = 1.051030612656967
$ws.Range("B22:F22").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.035953022230506
$arrIM[0,1] = 1.040657620605641
$arrIM[0,2] = 1.044534710338411
$arrIM[0,3] = 1.04597442427713
$arrIM[0,4] = 1.054682561882387
$ws.Range("I22:M22").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.034178443626861
$arrBF[0,2] = 1.041123567254314
$arrBF[0,3] = 1.042608669679068
$arrBF[0,4] = 1.051369119510191
$ws.Range("B23:F23").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036025244288313
$arrIM[0,1] = 1.040892320339576
$arrIM[0,2] = 1.044745091241385
$arrIM[0,3] = 1.046224636701035
$arrIM[0,4] = 1.054952659876683
$ws.Range("I23:M23").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.035623650027111
$arrBF[0,2] = 1.042223201192374
$arrBF[0,3] = 1.043866871164812
$arrBF[0,4] = 1.052703019460327
$ws.Range("B24:F24").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036305610619987
$arrIM[0,1] = 1.041815734061015
$arrIM[0,2] = 1.045571688450067
$arrIM[0,3] = 1.047209707502433
$arrIM[0,4] = 1.056015809537741
$ws.Range("I24:M24").Value2 = $arrIM

$arrBF = New-Object 'object[,]' 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.03730229971511
$arrBF[0,2] = 1.043498918422168
$arrBF[0,3] = 1.045329929011062
$arrBF[0,4] = 1.05425343493183
$ws.Range("B25:F25").Value2 = $arrBF

$arrIM = New-Object 'object[,]' 1,5
$arrIM[0,0] = 1.036622722686599
$arrIM[0,1] = 1.042886075342691
$arrIM[0,2] = 1.04652744715261
$arrIM[0,3] = 1.048352806574322
$arrIM[0,4] = 1.057249077680434
$ws.Range("I25:M25").Value2 = $arrIM
